$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 247.0944516666667
$ws.Range("N2").Value = 741.283355
$ws.Range("O2").Value = 0.8050739182622993
$ws.Range("P2").Value = 0.8050739182622993
$ws.Range("Q2").Value = 18.10790506630556
$ws.Range("R2").Value = 162.97114559675
$ws.Range("S2").Value = 0.8050739182622993
$ws.Range("T2").Value = 0.8050739182622993

# Row 3
$ws.Range("O3").Value = 0.1379009747488701
$ws.Range("P3").Value = 0.13790097474887
$ws.Range("S3").Value = 0.1379009747488701
$ws.Range("T3").Value = 0.13790097474887

# Row 4
$ws.Range("M4").Value = 11.590146
$ws.Range("N4").Value = 34.770438
$ws.Range("O4").Value = 0.03776258103132013
$ws.Range("P4").Value = 0.03776258103132013
$ws.Range("Q4").Value = 0.8493645326999998
$ws.Range("R4").Value = 7.644280794299999
$ws.Range("S4").Value = 0.03776258103132013
$ws.Range("T4").Value = 0.03776258103132013

# Row 5
$ws.Range("M5").Value = 5.912082333333333
$ws.Range("N5").Value = 17.736247
$ws.Range("O5").Value = 0.01926252595751047
$ws.Range("P5").Value = 0.01926252595751047
$ws.Range("Q5").Value = 0.4332571003277777
$ws.Range("R5").Value = 3.899313902949999
$ws.Range("S5").Value = 0.01926252595751047
$ws.Range("T5").Value = 0.01926252595751047
